$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Cells.Item(28,1).Value = 112551671
$ws.Cells.Item(28,2).Value = 78713
$ws.Cells.Item(28,3).Value = "Ovaliderad"
$ws.Cells.Item(28,4).Value = "NT"
$ws.Cells.Item(28,5).Value = 6458
$ws.Cells.Item(28,6).Value = "Lunglav"
$ws.Cells.Item(28,7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(28,8).Value = "(L.) Hoffm."
$ws.Cells.Item(28,16).Value = "Väster Naturreservatet Kullarna, Vrm"
$ws.Cells.Item(28,17).Value = 386308
$ws.Cells.Item(28,18).Value = 6713631
$ws.Cells.Item(28,19).Value = 10
$ws.Cells.Item(28,20).Value = "Värmland"
$ws.Cells.Item(28,21).Value = "Torsby"
$ws.Cells.Item(28,22).Value = "Värmland"
$ws.Cells.Item(28,23).Value = "Dalby"
$ws.Cells.Item(28,25).Value = "'2022-08-24"
$ws.Cells.Item(28,27).Value = "'2022-08-24"
$ws.Cells.Item(28,29).Value = "På asp."
$ws.Cells.Item(28,30).Value = $false
$ws.Cells.Item(28,31).Value = $false
$ws.Cells.Item(28,33).Value = $false
$ws.Cells.Item(28,49).Value = "Anders Boström"
$ws.Cells.Item(28,50).Value = "Anders Boström"

# Row 29
$ws.Cells.Item(29,1).Value = 112551677
$ws.Cells.Item(29,2).Value = 77650
$ws.Cells.Item(29,3).Value = "Ovaliderad"
$ws.Cells.Item(29,4).Value = "NT"
$ws.Cells.Item(29,5).Value = 6425
$ws.Cells.Item(29,6).Value = "Garnlav"
$ws.Cells.Item(29,7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(29,8).Value = "(Ach.) Ach."
$ws.Cells.Item(29,16).Value = "Väster Naturreservatet Kullarna, Vrm"
$ws.Cells.Item(29,17).Value = 386047
$ws.Cells.Item(29,18).Value = 6713439
$ws.Cells.Item(29,19).Value = 10
$ws.Cells.Item(29,20).Value = "Värmland"
$ws.Cells.Item(29,21).Value = "Torsby"
$ws.Cells.Item(29,22).Value = "Värmland"
$ws.Cells.Item(29,23).Value = "Dalby"
$ws.Cells.Item(29,25).Value = "'2022-08-24"
$ws.Cells.Item(29,27).Value = "'2022-08-24"
$ws.Cells.Item(29,30).Value = $false
$ws.Cells.Item(29,31).Value = $false
$ws.Cells.Item(29,33).Value = $false
$ws.Cells.Item(29,49).Value = "Anders Boström"
$ws.Cells.Item(29,50).Value = "Anders Boström"

# Row 30
$ws.Cells.Item(30,1).Value = 112551702
$ws.Cells.Item(30,2).Value = 78713
$ws.Cells.Item(30,3).Value = "Ovaliderad"
$ws.Cells.Item(30,4).Value = "NT"
$ws.Cells.Item(30,5).Value = 6458
$ws.Cells.Item(30,6).Value = "Lunglav"
$ws.Cells.Item(30,7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(30,8).Value = "(L.) Hoffm."
$ws.Cells.Item(30,16).Value = "Väster Naturreservatet Kullarna, Vrm"
$ws.Cells.Item(30,17).Value = 385983
$ws.Cells.Item(30,18).Value = 6713397
$ws.Cells.Item(30,19).Value = 10
$ws.Cells.Item(30,20).Value = "Värmland"
$ws.Cells.Item(30,21).Value = "Torsby"
$ws.Cells.Item(30,22).Value = "Värmland"
$ws.Cells.Item(30,23).Value = "Dalby"
$ws.Cells.Item(30,25).Value = "'2022-08-24"
$ws.Cells.Item(30,27).Value = "'2022-08-24"
$ws.Cells.Item(30,29).Value = "På asp."
$ws.Cells.Item(30,30).Value = $false
$ws.Cells.Item(30,31).Value = $false
$ws.Cells.Item(30,33).Value = $false
$ws.Cells.Item(30,49).Value = "Anders Boström"
$ws.Cells.Item(30,50).Value = "Anders Boström"

# Row 31
$ws.Cells.Item(31,1).Value = 112551681
$ws.Cells.Item(31,2).Value = 96735
$ws.Cells.Item(31,3).Value = "Ovaliderad"
$ws.Cells.Item(31,4).Value = "VU"
$ws.Cells.Item(31,5).Value = 220787
$ws.Cells.Item(31,6).Value = "Knärot"
$ws.Cells.Item(31,7).Value = "Goodyera repens"
$ws.Cells.Item(31,8).Value = "(L.) R. Br."
$ws.Cells.Item(31,16).Value = "Väster Naturreservatet Kullarna, Vrm"
$ws.Cells.Item(31,17).Value = 386295
$ws.Cells.Item(31,18).Value = 6713615
$ws.Cells.Item(31,19).Value = 10
$ws.Cells.Item(31,20).Value = "Värmland"
$ws.Cells.Item(31,21).Value = "Torsby"
$ws.Cells.Item(31,22).Value = "Värmland"
$ws.Cells.Item(31,23).Value = "Dalby"
$ws.Cells.Item(31,25).Value = "'2022-08-24"
$ws.Cells.Item(31,27).Value = "'2022-08-24"
$ws.Cells.Item(31,30).Value = $false
$ws.Cells.Item(31,31).Value = $false
$ws.Cells.Item(31,33).Value = $false
$ws.Cells.Item(31,49).Value = "Anders Boström"
$ws.Cells.Item(31,50).Value = "Anders Boström"

# Row 32
$ws.Cells.Item(32,1).Value = 112551679
$ws.Cells.Item(32,2).Value = 77650
$ws.Cells.Item(32,3).Value = "Ovaliderad"
$ws.Cells.Item(32,4).Value = "NT"
$ws.Cells.Item(32,5).Value = 6425
$ws.Cells.Item(32,6).Value = "Garnlav"
$ws.Cells.Item(32,7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(32,8).Value = "(Ach.) Ach."
$ws.Cells.Item(32,16).Value = "Väster Naturreservatet Kullarna, Vrm"
$ws.Cells.Item(32,17).Value = 386245
$ws.Cells.Item(32,18).Value = 6713552
$ws.Cells.Item(32,19).Value = 10
$ws.Cells.Item(32,20).Value = "Värmland"
$ws.Cells.Item(32,21).Value = "Torsby"
$ws.Cells.Item(32,22).Value = "Värmland"
$ws.Cells.Item(32,23).Value = "Dalby"
$ws.Cells.Item(32,25).Value = "'2022-08-24"
$ws.Cells.Item(32,27).Value = "'2022-08-24"
$ws.Cells.Item(32,30).Value = $false
$ws.Cells.Item(32,31).Value = $false
$ws.Cells.Item(32,33).Value = $false
$ws.Cells.Item(32,49).Value = "Anders Boström"
$ws.Cells.Item(32,50).Value = "Anders Boström"
